$d = $word.ActiveDocument
$bullet = [char]0x2022

# The "KEY ACHIEVEMENTS AND IMPACT" section contains bullet paragraphs whose
# text also shows up (duplicated / near-duplicated) elsewhere in the document
# (e.g. under "PROFESSIONAL EXPERIENCE"), so this edit is targeted by first
# locating the "Impact" sub-heading that precedes these specific bullets
# rather than doing a blind document-wide Find/Replace.

$targetHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $targetHeading = $i
        break
    }
}

if ($targetHeading -eq $null) {
    throw "Could not locate 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# "Impact" sub-heading is the paragraph right after the section heading.
$impactIdx = $targetHeading + 1
$impactPara = $d.Paragraphs.Item($impactIdx)
if ($impactPara.Range.Text.TrimEnd([char]13) -ne "Impact") {
    throw "Unexpected paragraph after KEY ACHIEVEMENTS AND IMPACT heading"
}

# The six bullet paragraphs follow directly after "Impact". Paragraphs.Item()
# resolves by live index, not a stable handle, so re-fetch by index after any
# edit that changes the paragraph count (e.g. Range.Delete of a whole
# paragraph).
$bulletBase = $impactIdx

function Get-Bullet($offset) {
    return $d.Paragraphs.Item($bulletBase + $offset)
}

# Sanity-check the bullets are the ones we expect before mutating anything.
if ((Get-Bullet 1).Range.Text -notmatch "Built redistricting platform used by thousands") {
    throw "Bullet 1 text mismatch"
}
if ((Get-Bullet 2).Range.Text -notmatch "Designed ETL pipelines using PySpark") {
    throw "Bullet 2 text mismatch"
}
if ((Get-Bullet 3).Range.Text -notmatch "Trigonometric algorithm for boundary estimation") {
    throw "Bullet 3 text mismatch"
}
if ((Get-Bullet 4).Range.Text -notmatch "Discovered systematic race coding errors") {
    throw "Bullet 4 text mismatch"
}
if ((Get-Bullet 5).Range.Text -notmatch "Achieved 87% prediction accuracy") {
    throw "Bullet 5 text mismatch"
}
if ((Get-Bullet 6).Range.Text -notmatch "Built cloud-based data warehouse") {
    throw "Bullet 6 text mismatch"
}

# Rewrite bullets 1-3 in place first (their indices do not shift). Note: do
# NOT append a trailing `r` — Range.Text already replaces up to (but not
# including) the existing paragraph mark, so adding `r` would insert an
# extra empty paragraph.
(Get-Bullet 1).Range.Text = "$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
(Get-Bullet 2).Range.Text = "$bullet `$4.7M savings enabled nonprofit access"
(Get-Bullet 3).Range.Text = "$bullet Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

# Remove bullets 4 and 5 entirely (delete the whole paragraph, mark included).
# Deleting the later one (5) first keeps bullet 4's index stable for its own
# delete.
(Get-Bullet 5).Range.Delete()
(Get-Bullet 4).Range.Delete()

# Bullet 6 has now shifted up to offset 4 (two paragraphs were removed).
(Get-Bullet 4).Range.Text = "$bullet Real-time collaboration at national scale"
